# "Add files via upload" - the single data record on Sheet1 was replaced
# with a new hospital entry. The header row (row 1) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hospital4042"
$ws.Range("B2").Value = "hospital4042"
$ws.Range("C2").Value = "Sub District Hospital"
$ws.Range("D2").Value = "0832 231 4042"
$ws.Range("E2").Value = "Ponda"
$ws.Range("F2").Value = "Rajesh Naik"
# Leading apostrophe keeps this all-digit phone number stored as text
# instead of being auto-converted into a number.
$ws.Range("G2").Value = "'8450124593"

# Leave the whole second row selected, as it was after the upload.
$ws.Rows("2:2").Select()
